# Updates the two-digit / one-digit division worksheet numbers.
# The document has a single table; the populated rows are 1, 5, 9, 13, 17
# (1-based), each with 5 columns of "NN÷N=" style problems. We replace the
# text of each populated cell directly by (row, column) position so that
# duplicate source strings (e.g. "56÷3=" appearing twice) are handled
# correctly and unambiguously.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$replacements = @(
    @{ Row = 1;  Col = 1; Old = "88÷5="; New = "64÷3=" },
    @{ Row = 1;  Col = 2; Old = "37÷2="; New = "67÷8=" },
    @{ Row = 1;  Col = 3; Old = "88÷4="; New = "75÷7=" },
    @{ Row = 1;  Col = 4; Old = "56÷3="; New = "42÷4=" },
    @{ Row = 1;  Col = 5; Old = "56÷3="; New = "53÷9=" },

    @{ Row = 5;  Col = 1; Old = "70÷4="; New = "36÷4=" },
    @{ Row = 5;  Col = 2; Old = "48÷6="; New = "25÷6=" },
    @{ Row = 5;  Col = 3; Old = "93÷5="; New = "76÷3=" },
    @{ Row = 5;  Col = 4; Old = "62÷9="; New = "22÷8=" },
    @{ Row = 5;  Col = 5; Old = "14÷6="; New = "10÷9=" },

    @{ Row = 9;  Col = 1; Old = "79÷8="; New = "51÷9=" },
    @{ Row = 9;  Col = 2; Old = "27÷8="; New = "13÷2=" },
    @{ Row = 9;  Col = 3; Old = "26÷2="; New = "37÷2=" },
    @{ Row = 9;  Col = 4; Old = "36÷3="; New = "96÷9=" },
    @{ Row = 9;  Col = 5; Old = "47÷3="; New = "76÷2=" },

    @{ Row = 13; Col = 1; Old = "94÷5="; New = "99÷3=" },
    @{ Row = 13; Col = 2; Old = "20÷9="; New = "25÷9=" },
    @{ Row = 13; Col = 3; Old = "91÷4="; New = "29÷6=" },
    @{ Row = 13; Col = 4; Old = "40÷9="; New = "96÷4=" },
    @{ Row = 13; Col = 5; Old = "96÷9="; New = "90÷6=" },

    @{ Row = 17; Col = 1; Old = "25÷5="; New = "87÷2=" },
    @{ Row = 17; Col = 2; Old = "35÷3="; New = "57÷9=" },
    @{ Row = 17; Col = 3; Old = "61÷9="; New = "42÷4=" },
    @{ Row = 17; Col = 4; Old = "48÷5="; New = "57÷7=" },
    @{ Row = 17; Col = 5; Old = "88÷3="; New = "49÷3=" }
)

foreach ($r in $replacements) {
    $cell = $t.Cell($r.Row, $r.Col)
    $range = $cell.Range
    # Trim the trailing end-of-cell marker so only the visible text is
    # replaced; Find/ReplaceAll on a sub-range isn't reliably scoped when
    # the same text occurs elsewhere in the document (duplicate problems
    # such as "56÷3=" appear twice), so we assign the text directly.
    $range.MoveEnd(1, -1)  # wdCharacter = 1; trim the end-of-cell mark
    $range.Text = $r.New
}
